# Value of a Statistical Life.xlsx - update from China (World Bank) VoaSL
# to US EPA VoaSL, dropping the intermediate "Data" worksheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "About" sheet: new source (US EPA instead of World Bank),
#    new notes text, and new currency-year-adjustment block (rows 15-18).
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("B3").Value = "U.S. Environmental Protection Agency"
$about.Range("B4").Value = 2013
$about.Range("B5").Value = "Frequently Asked Questions on Mortality Risk Valuation"
$about.Range("B6").Value = "http://yosemite.epa.gov/EE%5Cepa%5Ceed.nsf/webpages/MortalityRiskValuation.html#whatvalue"
$about.Range("B7").Value = """What value of statistical life does EPA use?"""

# remove the hyperlink that used to point at the World Bank page - the new
# text in B6 is plain text, not a live hyperlink
$about.Hyperlinks.Delete()

# Notes section (A9:A13) keeps the same wording, only row 9's caption
# ("Notes") stays; rows 10-13 keep their original text.

# New "Currency Year Adjustment" block (same bold look as the "Notes" caption)
$about.Range("A15").Font.Bold = $true
$about.Range("A15").Value = "Currency Year Adjustment"
$about.Range("A16").Value = "We adjust 2006 dollars to 2012 dollars using the following conversion factor:"
$about.Range("A17").Value = 1.141
$about.Range("A18").Value = "See ""cpi.xlsx"" in the InputData folder for source information."

# ---------------------------------------------------------------------
# 2. Update the "VoaSL" sheet: new headers/labels and a new formula that
#    references About!A17 instead of the (to-be-removed) Data sheet.
# ---------------------------------------------------------------------
$voasl = $wb.Worksheets.Item("VoaSL")

$voasl.Range("A1").Value = "2012$/life"
$voasl.Range("A2").Value = "Value"

$voasl.Range("B2").ClearFormats()
$voasl.Range("B2").Formula = "=7.4*10^6*About!A17"

# ---------------------------------------------------------------------
# 3. Delete the intermediate "Data" worksheet - its one formula has been
#    folded directly into VoaSL!B2 above.
# ---------------------------------------------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Data").Delete()
$excel.DisplayAlerts = $true

# ---------------------------------------------------------------------
# 4. Restore a sane view state: VoaSL selected at A3, but "About" is the
#    tab that is actually active/selected when the workbook is opened.
# ---------------------------------------------------------------------
$voasl.Range("A3").Select()
$about.Activate()
$about.Range("A1").Select()
